# Weekly update: insert a new record at row 13 (new market-day observation),
# pushing the existing rows 13-41 down to 14-42.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 13; this shifts rows 13..41 down to 14..42
# and keeps their existing values/formatting intact.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new observation.
$ws.Cells.Item(13, 1).Value = 10
$ws.Cells.Item(13, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(13, 3).Value = "La Araucanía"
$ws.Cells.Item(13, 4).Value = 44497
$ws.Cells.Item(13, 5).Value = 9
$ws.Cells.Item(13, 6).Value = 100112026
$ws.Cells.Item(13, 7).Value = "Haba"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 200
$ws.Cells.Item(13, 11).Value = 8000
$ws.Cells.Item(13, 12).Value = 9000
$ws.Cells.Item(13, 13).Value = 8500
$ws.Cells.Item(13, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(13, 15).Value = "Región Metropolitana"
$ws.Cells.Item(13, 16).Value = 340
$ws.Cells.Item(13, 17).Value = 25
$ws.Cells.Item(13, 18).Value = "Hortaliza"
